$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics now that trade #10 has closed
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.99   # Current Capital
$wsSummary.Range("B4").Value = -0.01     # Total P&L $
$wsSummary.Range("B5").Value = -0.02     # Total P&L %
$wsSummary.Range("B6").Value = 10        # Total Trades
$wsSummary.Range("B8").Value = 5         # Losing Trades
$wsSummary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.98999999999999  # Capital
$wsStatus.Range("D4").Value = 10                 # Trades
$wsStatus.Range("E4").Value = -0.01               # P&L $
$wsStatus.Range("F4").Value = -0.01               # P&L %
$wsStatus.Range("G4").Value = 40                  # Win Rate %

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append new row 11 for trade #10
# ---------------------------------------------------------------------------
function Add-Trade10Row($ws) {
    $ws.Cells.Item(11, 1).Value = 10             # Trade #
    $ws.Cells.Item(11, 2).NumberFormat = "@"
    $ws.Cells.Item(11, 2).Value = "2026-02-17"   # Date (keep as text)
    $ws.Cells.Item(11, 2).Style = "Normal"
    $ws.Cells.Item(11, 3).Value = "07:58:19"     # Time
    $ws.Cells.Item(11, 4).Value = "MarketMaking" # Strategy
    $ws.Cells.Item(11, 5).Value = "UP"           # Side
    $ws.Cells.Item(11, 6).Value = 0.11           # Entry Price
    $ws.Cells.Item(11, 7).Value = 0.09           # Exit Price
    $ws.Cells.Item(11, 8).Value = "CLOSED"       # Status
    $ws.Cells.Item(11, 9).Value = -18.1818       # P&L %
    $ws.Cells.Item(11, 10).Value = -0.02         # P&L $
    $ws.Cells.Item(11, 11).Value = 99.98999999999999  # Capital After
    $ws.Cells.Item(11, 12).Value = 0             # Entry Slippage (bps)
    $ws.Cells.Item(11, 13).Value = 0             # Exit Slippage (bps)
    $ws.Cells.Item(11, 14).Value = 0.6           # Confidence
    $ws.Cells.Item(11, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(11, 16).Value = "early_exit"  # Exit Reason
    $ws.Cells.Item(11, 17).Value = 0.13          # Duration (min)
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade10Row $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade10Row $wsMarketMaking
